$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 2.191602
$ws.Cells.Item(2, 8).Value = 6.574806000000001
$ws.Cells.Item(2, 9).Value = 0.07674610985252207
$ws.Cells.Item(2, 10).Value = 0.07674610985252209
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 218.7785543333333
$ws.Cells.Item(2, 14).Value = 656.3356630000001
$ws.Cells.Item(2, 15).Value = 0.7837094150017259
$ws.Cells.Item(2, 16).Value = 0.7837094150017259
$ws.Cells.Item(2, 17).Value = 479.4755172340421
$ws.Cells.Item(2, 18).Value = 4315.279655106378
$ws.Cells.Item(2, 19).Value = 0.06014664885617826
$ws.Cells.Item(2, 20).Value = 0.06014664885617828

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 2.191602
$ws.Cells.Item(3, 8).Value = 6.574806000000001
$ws.Cells.Item(3, 9).Value = 0.07674610985252207
$ws.Cells.Item(3, 10).Value = 0.07674610985252209
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 46.29469433333333
$ws.Cells.Item(3, 14).Value = 138.884083
$ws.Cells.Item(3, 15).Value = 0.1658370397602197
$ws.Cells.Item(3, 16).Value = 0.1658370397602197
$ws.Cells.Item(3, 17).Value = 101.459544690322
$ws.Cells.Item(3, 18).Value = 913.1359022128981
$ws.Cells.Item(3, 19).Value = 0.01272734767105489
$ws.Cells.Item(3, 20).Value = 0.01272734767105489

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 2.191602
$ws.Cells.Item(4, 8).Value = 6.574806000000001
$ws.Cells.Item(4, 9).Value = 0.07674610985252207
$ws.Cells.Item(4, 10).Value = 0.07674610985252209
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 6.978882
$ws.Cells.Item(4, 14).Value = 20.936646
$ws.Cells.Item(4, 15).Value = 0.02499977909741928
$ws.Cells.Item(4, 16).Value = 0.02499977909741927
$ws.Cells.Item(4, 17).Value = 15.294931748964
$ws.Cells.Item(4, 18).Value = 137.654385740676
$ws.Cells.Item(4, 19).Value = 0.001918635792899325
$ws.Cells.Item(4, 20).Value = 0.001918635792899325

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 2.191602
$ws.Cells.Item(5, 8).Value = 6.574806000000001
$ws.Cells.Item(5, 9).Value = 0.07674610985252207
$ws.Cells.Item(5, 10).Value = 0.07674610985252209
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 7.105616
$ws.Cells.Item(5, 14).Value = 21.316848
$ws.Cells.Item(5, 15).Value = 0.02545376614063513
$ws.Cells.Item(5, 16).Value = 0.02545376614063513
$ws.Cells.Item(5, 17).Value = 15.572682236832
$ws.Cells.Item(5, 18).Value = 140.154140131488
$ws.Cells.Item(5, 19).Value = 0.001953477532389591
$ws.Cells.Item(5, 20).Value = 0.001953477532389591

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 16.59481266666667
$ws.Cells.Item(6, 8).Value = 49.78443799999999
$ws.Cells.Item(6, 9).Value = 0.5811216251390648
$ws.Cells.Item(6, 10).Value = 0.5811216251390647
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 218.7785543333333
$ws.Cells.Item(6, 14).Value = 656.3356630000001
$ws.Cells.Item(6, 15).Value = 0.7837094150017259
$ws.Cells.Item(6, 16).Value = 0.7837094150017259
$ws.Cells.Item(6, 17).Value = 3630.589124645822
$ws.Cells.Item(6, 18).Value = 32675.30212181239
$ws.Cells.Item(6, 19).Value = 0.4554304888825887
$ws.Cells.Item(6, 20).Value = 0.4554304888825886

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 16.59481266666667
$ws.Cells.Item(7, 8).Value = 49.78443799999999
$ws.Cells.Item(7, 9).Value = 0.5811216251390648
$ws.Cells.Item(7, 10).Value = 0.5811216251390647
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 46.29469433333333
$ws.Cells.Item(7, 14).Value = 138.884083
$ws.Cells.Item(7, 15).Value = 0.1658370397602197
$ws.Cells.Item(7, 16).Value = 0.1658370397602197
$ws.Cells.Item(7, 17).Value = 768.2517799222616
$ws.Cells.Item(7, 18).Value = 6914.266019300353
$ws.Cells.Item(7, 19).Value = 0.09637149005371055
$ws.Cells.Item(7, 20).Value = 0.09637149005371054

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 16.59481266666667
$ws.Cells.Item(8, 8).Value = 49.78443799999999
$ws.Cells.Item(8, 9).Value = 0.5811216251390648
$ws.Cells.Item(8, 10).Value = 0.5811216251390647
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 6.978882
$ws.Cells.Item(8, 14).Value = 20.936646
$ws.Cells.Item(8, 15).Value = 0.02499977909741928
$ws.Cells.Item(8, 16).Value = 0.02499977909741927
$ws.Cells.Item(8, 17).Value = 115.813239412772
$ws.Cells.Item(8, 18).Value = 1042.319154714948
$ws.Cells.Item(8, 19).Value = 0.01452791225720991
$ws.Cells.Item(8, 20).Value = 0.01452791225720991

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 16.59481266666667
$ws.Cells.Item(9, 8).Value = 49.78443799999999
$ws.Cells.Item(9, 9).Value = 0.5811216251390648
$ws.Cells.Item(9, 10).Value = 0.5811216251390647
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 7.105616
$ws.Cells.Item(9, 14).Value = 21.316848
$ws.Cells.Item(9, 15).Value = 0.02545376614063513
$ws.Cells.Item(9, 16).Value = 0.02545376614063513
$ws.Cells.Item(9, 17).Value = 117.9163664012693
$ws.Cells.Item(9, 18).Value = 1061.247297611424
$ws.Cells.Item(9, 19).Value = 0.01479173394555559
$ws.Cells.Item(9, 20).Value = 0.01479173394555558

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 9.770107666666666
$ws.Cells.Item(10, 8).Value = 29.310323
$ws.Cells.Item(10, 9).Value = 0.342132265008413
$ws.Cells.Item(10, 10).Value = 0.342132265008413
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 218.7785543333333
$ws.Cells.Item(10, 14).Value = 656.3356630000001
$ws.Cells.Item(10, 15).Value = 0.7837094150017259
$ws.Cells.Item(10, 16).Value = 0.7837094150017259
$ws.Cells.Item(10, 17).Value = 2137.49003099435
$ws.Cells.Item(10, 18).Value = 19237.41027894915
$ws.Cells.Item(10, 19).Value = 0.2681322772629589
$ws.Cells.Item(10, 20).Value = 0.2681322772629589

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 9.770107666666666
$ws.Cells.Item(11, 8).Value = 29.310323
$ws.Cells.Item(11, 9).Value = 0.342132265008413
$ws.Cells.Item(11, 10).Value = 0.342132265008413
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 46.29469433333333
$ws.Cells.Item(11, 14).Value = 138.884083
$ws.Cells.Item(11, 15).Value = 0.1658370397602197
$ws.Cells.Item(11, 16).Value = 0.1658370397602197
$ws.Cells.Item(11, 17).Value = 452.3041480320899
$ws.Cells.Item(11, 18).Value = 4070.737332288809
$ws.Cells.Item(11, 19).Value = 0.0567382020354542
$ws.Cells.Item(11, 20).Value = 0.0567382020354542

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 9.770107666666666
$ws.Cells.Item(12, 8).Value = 29.310323
$ws.Cells.Item(12, 9).Value = 0.342132265008413
$ws.Cells.Item(12, 10).Value = 0.342132265008413
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 6.978882
$ws.Cells.Item(12, 14).Value = 20.936646
$ws.Cells.Item(12, 15).Value = 0.02499977909741928
$ws.Cells.Item(12, 16).Value = 0.02499977909741927
$ws.Cells.Item(12, 17).Value = 68.18442853296199
$ws.Cells.Item(12, 18).Value = 613.6598567966579
$ws.Cells.Item(12, 19).Value = 0.008553231047310038
$ws.Cells.Item(12, 20).Value = 0.008553231047310036

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 9.770107666666666
$ws.Cells.Item(13, 8).Value = 29.310323
$ws.Cells.Item(13, 9).Value = 0.342132265008413
$ws.Cells.Item(13, 10).Value = 0.342132265008413
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 7.105616
$ws.Cells.Item(13, 14).Value = 21.316848
$ws.Cells.Item(13, 15).Value = 0.02545376614063513
$ws.Cells.Item(13, 16).Value = 0.02545376614063513
$ws.Cells.Item(13, 17).Value = 69.42263335798934
$ws.Cells.Item(13, 18).Value = 624.8037002219039
$ws.Cells.Item(13, 19).Value = 0.00870855466268995
$ws.Cells.Item(13, 20).Value = 0.008708554662689948
